$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Personal Summary paragraph: update wording + merge with the blank
#    paragraph that used to follow it (the diff removes the paragraph break
#    that separated the summary text from the next run).
# ---------------------------------------------------------------------------
$oldSummary = "UI/Front-end developer with proven track-record over five years in creating clean, accessible and user-friendly responsive websites. Self motivated and adaptable, independent minded, with a keen eye for detail. Experienced in working with Hybris CMS and Agile methodologies."
$newSummary = "UI/Front-end developer with proven track-record over five years in creating accessible and user-friendly responsive websites. Self motivated and independent minded, with a keen eye for detail. Two months of commercial use of AngularJS on B2B project. Experienced in working with Hybris CMS and Agile methodologies."

$d.Content.Find.Execute($oldSummary, $false, $false, $false, $false, $false, `
                         $true, 1, $false, $newSummary, 2)

# Locate the paragraph that now holds the updated summary text (by index, to
# avoid relying on `.Next`, which can hand back a range whose `.Text` is
# unreliable to compare against) and merge it with the following paragraph
# when that paragraph is blank, by deleting the paragraph mark between them.
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -like "*Two months of commercial use of AngularJS on B2B project*") {
        $targetIndex = $i
    }
}

if ($targetIndex -gt 0) {
    $nextIndex = $targetIndex + 1
    $nextPara = $d.Paragraphs($nextIndex)
    $nextText = $nextPara.Range.Text
    $trimmed = $nextText.Trim()
    if ($trimmed.Length -eq 0) {
        $summaryPara = $d.Paragraphs($targetIndex)
        $mark = $d.Range($summaryPara.Range.End - 1, $summaryPara.Range.End)
        $mark.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2) "Developing and integrating front-end..." bullet: Grunt.js -> AngularJS
# ---------------------------------------------------------------------------
$oldBullet = "Developing and integrating front-end with Hybris[y] CMS for number of e-commerce websites using HTML5, CSS3, Less/SASS, Bootstrap 3, jQuery and Grunt.js"
$newBullet = "Developing and integrating front-end with Hybris[y] CMS for number of e-commerce websites using HTML5, CSS3, Less/SASS, Bootstrap 3, jQuery and AngularJS"

$d.Content.Find.Execute($oldBullet, $false, $false, $false, $false, $false, `
                         $true, 1, $false, $newBullet, 2)
